$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C holds the "Förändrad" (last changed) date; bump it from
# 45181 (2023-09-12) to 45182 (2023-09-13) for every data row.
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45182
